$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data to insert right after the header row (row 1),
# pushing the existing data down by 8 rows.
$newTopRows = @(
    @(-2.872348141670227, 9.113740730285643, -1.157746517658234),
    @(-2.973343849182129, 9.098371124267578, -1.108406949043274),
    @(-3.070270323753357, 9.036239910125733, -1.239402884244919),
    @(-3.125558137893677, 8.97771692276001,  -1.272288262844086),
    @(-3.148761796951294, 8.941205215454101, -1.218037897348404),
    @(-3.052603721618652, 8.98315830230713,  -1.222838401794433),
    @(-3.015902495384216, 9.072346258163453, -1.212291812896728),
    @(-3.051667261123657, 9.041837882995607, -1.251752722263336)
)

# New rows of data to append at the bottom of the existing data.
$newBottomRows = @(
    @(2.46818006038666,  7.359129667282104, -4.541467189788817),
    @(3.081884574890136, 7.212435054779053, -4.413347625732421)
)

$insertCount = $newTopRows.Length

# Insert the needed number of blank rows right after the header (before row 2),
# shifting all existing data down.
$insertRange = $ws.Range("A2:C" + (1 + $insertCount))
$insertRange.Insert()

# Fill in the newly inserted rows with their values.
for ($i = 0; $i -lt $newTopRows.Length; $i++) {
    $r = 2 + $i
    $row = $newTopRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}

# Append the two new rows at the end (after what is now row 29).
$lastRow = 1 + $insertCount + 20
for ($i = 0; $i -lt $newBottomRows.Length; $i++) {
    $r = $lastRow + 1 + $i
    $row = $newBottomRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}
